$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stock List")

$ws.Range("B2:C2").Value = "GROWW"
$ws.Range("D2").Value = 162.66
$ws.Range("E2").Value = 9.513199999999999
$ws.Range("H2").Value = 91696.4308
$ws.Range("B3:C3").Value = "TMCV"
$ws.Range("D3").Value = 320.05
$ws.Range("E3").Value = 0.7714
$ws.Range("H3").Value = 116950.8444
$ws.Range("B4:C4").Value = "SMLMAH"
$ws.Range("D4").Value = 2870
$ws.Range("E4").Value = 0.4128
$ws.Range("H4").Value = 0
$ws.Range("B5:C5").Value = "IBULLSLTD"
$ws.Range("D5").Value = 21.02
$ws.Range("E5").Value = "N/A"
$ws.Range("H5").Value = 0
$ws.Range("B6:C6").Value = "CHOICEGOLD"
$ws.Range("D6").Value = 121.4
$ws.Range("E6").Value = -2.6385
$ws.Range("H6").Value = 0
$ws.Range("B7:C7").Value = "LENSKART"
$ws.Range("D7").Value = 420
$ws.Range("E7").Value = 2.6769
$ws.Range("H7").Value = 70964.86500000001
$ws.Range("B8:C8").Value = "STUDDS"
$ws.Range("D8").Value = 548
$ws.Range("E8").Value = -0.8235
$ws.Range("H8").Value = 2174.4721
$ws.Range("B9:C9").Value = "PIRAMALFIN"
$ws.Range("D9").Value = 1564.7
$ws.Range("E9").Value = 4.9993
$ws.Range("H9").Value = 0
$ws.Range("B10:C10").Value = "ORKLAINDIA"
$ws.Range("D10").Value = 679.9
$ws.Range("E10").Value = 1.9111
$ws.Range("H10").Value = 9139.236500000001
$ws.Range("B11:C11").Value = "GROWWSC250"
$ws.Range("D11").Value = 9.890000000000001
$ws.Range("E11").Value = 0.5081
$ws.Range("H11").Value = 0
$ws.Range("B12:C12").Value = "MIDWESTLTD"
$ws.Range("D12").Value = 1284.3
$ws.Range("E12").Value = 9.544499999999999
$ws.Range("H12").Value = 4239.4941
$ws.Range("B13:C13").Value = "NIFTYCASE"
$ws.Range("D13").Value = 10.32
$ws.Range("E13").Value = 0.4868
$ws.Range("H13").Value = 0
$ws.Range("B14:C14").Value = "MOMENTUM30"
$ws.Range("D14").Value = 31.91
$ws.Range("E14").Value = 0.2828
$ws.Range("H14").Value = 0
$ws.Range("B15:C15").Value = "CANHLIFE"
$ws.Range("D15").Value = 120.06
$ws.Range("E15").Value = -0.2078
$ws.Range("H15").Value = 11429.45
$ws.Range("B16:C16").Value = "FLEXIADD"
$ws.Range("D16").Value = 10.89
$ws.Range("E16").Value = 0.554
$ws.Range("H16").Value = 0
$ws.Range("B17:C17").Value = "MOENERGY"
$ws.Range("D17").Value = 36.51
$ws.Range("E17").Value = 0.3573
$ws.Range("H17").Value = 0
$ws.Range("B18:C18").Value = "MONIFTY100"
$ws.Range("D18").Value = 26.7
$ws.Range("E18").Value = 0.1125
$ws.Range("H18").Value = 0
$ws.Range("B19:C19").Value = "RUBICON"
$ws.Range("D19").Value = 693
$ws.Range("E19").Value = -2.914
$ws.Range("H19").Value = 11759.8771
$ws.Range("B20:C20").Value = "CRAMC"
$ws.Range("D20").Value = 295.95
$ws.Range("E20").Value = 1.318
$ws.Range("H20").Value = 87.7936
$ws.Range("B21:C21").Value = "LGEINDIA"
$ws.Range("D21").Value = 1625
$ws.Range("E21").Value = 0.445
$ws.Range("H21").Value = 109811.7976
$ws.Range("B22:C22").Value = "TATACAP"
$ws.Range("D22").Value = 325.05
$ws.Range("E22").Value = -0.1383
$ws.Range("H22").Value = 138170.4872
$ws.Range("B23:C23").Value = "WEWORK"
$ws.Range("D23").Value = 613.85
$ws.Range("E23").Value = -1.2627
$ws.Range("H23").Value = 8332.226000000001
$ws.Range("B24:C24").Value = "GROWWRLTY"
$ws.Range("D24").Value = 10.6
$ws.Range("E24").Value = 0.6647999999999999
$ws.Range("H24").Value = 0
$ws.Range("B25:C25").Value = "ADVANCE"
$ws.Range("D25").Value = 145.66
$ws.Range("E25").Value = 2.7294
$ws.Range("H25").Value = 911.5072
$ws.Range("B26:C26").Value = "OMFREIGHT"
$ws.Range("D26").Value = 90.62
$ws.Range("E26").Value = -0.2532
$ws.Range("H26").Value = 305.9414
$ws.Range("B27:C27").Value = "GLOTTIS"
$ws.Range("D27").Value = 70.90000000000001
$ws.Range("E27").Value = -3.8514
$ws.Range("H27").Value = 681.3797
$ws.Range("B28:C28").Value = "FABTECH"
$ws.Range("D28").Value = 251.21
$ws.Range("E28").Value = 1.919
$ws.Range("H28").Value = 1095.6232
$ws.Range("B29:C29").Value = "PACEDIGITK"
$ws.Range("D29").Value = 210.64
$ws.Range("E29").Value = -3.9664
$ws.Range("H29").Value = 4734.483
$ws.Range("B30:C30").Value = "JAINREC"
$ws.Range("D30").Value = 443.45
$ws.Range("E30").Value = 0.8529
$ws.Range("H30").Value = 15173.4232
$ws.Range("B31:C31").Value = "EPACKPEB"
$ws.Range("D31").Value = 330.3
$ws.Range("E31").Value = 0.8396
$ws.Range("H31").Value = 3290.3052
$ws.Range("B32:C32").Value = "BMWVENTLTD"
$ws.Range("D32").Value = 69.20999999999999
$ws.Range("E32").Value = -0.703
$ws.Range("H32").Value = 604.4036
$ws.Range("B33:C33").Value = "STYL"
$ws.Range("D33").Value = 324.85
$ws.Range("E33").Value = -2.0651
$ws.Range("H33").Value = 5367.0992
$ws.Range("B34:C34").Value = "JARO"
$ws.Range("D34").Value = 691.7
$ws.Range("E34").Value = 10.3014
$ws.Range("H34").Value = 1389.4209
$ws.Range("B35:C35").Value = "SOLARWORLD"
$ws.Range("D35").Value = 299.75
$ws.Range("E35").Value = -0.1
$ws.Range("H35").Value = 2600.613
$ws.Range("B36:C36").Value = "ARSSBL"
$ws.Range("D36").Value = 742
$ws.Range("E36").Value = 1.1726
$ws.Range("H36").Value = 4600.2698
$ws.Range("B37:C37").Value = "GANESHCP"
$ws.Range("D37").Value = 284
$ws.Range("E37").Value = -1.7131
$ws.Range("H37").Value = 1167.7321
$ws.Range("B38:C38").Value = "ATLANTAELE"
$ws.Range("D38").Value = 962.1
$ws.Range("E38").Value = -0.3986
$ws.Range("H38").Value = 7427.8295
$ws.Range("B39:C39").Value = "GKENERGY"
$ws.Range("D39").Value = 201.79
$ws.Range("E39").Value = 1.683
$ws.Range("H39").Value = 4024.9086
$ws.Range("B40:C40").Value = "SAATVIKGL"
$ws.Range("D40").Value = 455
$ws.Range("E40").Value = 0.6192
$ws.Range("H40").Value = 5747.6883
$ws.Range("B41:C41").Value = "IVALUE"
$ws.Range("D41").Value = 313.2
$ws.Range("E41").Value = 0.2561
$ws.Range("H41").Value = 1672.5859
$ws.Range("B42:C42").Value = "VMSTMT"
$ws.Range("D42").Value = 64.86
$ws.Range("E42").Value = -3.5682
$ws.Range("H42").Value = 333.8195
$ws.Range("B43:C43").Value = "EUROPRATIK"
$ws.Range("D43").Value = 371.55
$ws.Range("E43").Value = 5.1358
$ws.Range("H43").Value = 3611.748
$ws.Range("B44:C44").Value = "SHRINGARMS"
$ws.Range("D44").Value = 214.3
$ws.Range("E44").Value = -0.7319
$ws.Range("H44").Value = 2081.7757
$ws.Range("B45:C45").Value = "DEVX"
$ws.Range("D45").Value = 42.13
$ws.Range("E45").Value = -1.1961
$ws.Range("H45").Value = 384.5596
$ws.Range("B46:C46").Value = "URBANCO"
$ws.Range("D46").Value = 142.13
$ws.Range("E46").Value = 1.0164
$ws.Range("H46").Value = 20203.1402
$ws.Range("B47:C47").Value = "SML100CASE"
$ws.Range("D47").Value = 10.35
$ws.Range("E47").Value = 0.779
$ws.Range("H47").Value = 0
$ws.Range("B48:C48").Value = "AONEGOLD"
$ws.Range("D48").Value = 11.45
$ws.Range("E48").Value = -2.1368
$ws.Range("H48").Value = 0
$ws.Range("B49:C49").Value = "ELM250"
$ws.Range("D49").Value = 16.88
$ws.Range("E49").Value = 0.1186
$ws.Range("H49").Value = 0
$ws.Range("B50:C50").Value = "AMANTA"
$ws.Range("D50").Value = 122.54
$ws.Range("E50").Value = 2.2701
$ws.Range("H50").Value = 465.2533
$ws.Range("B51:C51").Value = "CPEDU"
$ws.Range("D51").Value = 287.95
$ws.Range("E51").Value = -0.501
$ws.Range("H51").Value = 526.5037
$ws.Range("B52:C52").Value = "AHCL"
$ws.Range("D52").Value = 144
$ws.Range("E52").Value = 0.7134
$ws.Range("H52").Value = 759.9601
$ws.Range("B53:C53").Value = "STLNETWORK"
$ws.Range("D53").Value = 23.69
$ws.Range("E53").Value = 0.68
$ws.Range("H53").Value = 1148.0783
$ws.Range("B54:C54").Value = "VIKRAN"
$ws.Range("D54").Value = 110.31
$ws.Range("E54").Value = 1.2297
$ws.Range("H54").Value = 2810.4554
$ws.Range("B55:C55").Value = "MANUFGBEES"
$ws.Range("D55").Value = 153.75
$ws.Range("E55").Value = -0.0325
$ws.Range("H55").Value = 0
$ws.Range("B56:C56").Value = "MEIL"
$ws.Range("D56").Value = 448.55
$ws.Range("E56").Value = 0.7072000000000001
$ws.Range("H56").Value = 1230.6457
$ws.Range("B57:C57").Value = "GROWWNXT50"
$ws.Range("D57").Value = 70.17
$ws.Range("E57").Value = 0.3002
$ws.Range("H57").Value = 0
$ws.Range("B58:C58").Value = "SHREEJISPG"
$ws.Range("D58").Value = 306.95
$ws.Range("E58").Value = 3.0379
$ws.Range("H58").Value = 4853.3348
$ws.Range("B59:C59").Value = "GEMAROMA"
$ws.Range("D59").Value = 186.02
$ws.Range("E59").Value = -0.8369
$ws.Range("H59").Value = 979.9165
$ws.Range("B60:C60").Value = "PATELRMART"
$ws.Range("D60").Value = 223
$ws.Range("E60").Value = 0.4685
$ws.Range("H60").Value = 741.3581
$ws.Range("B61:C61").Value = "VIKRAMSOLR"
$ws.Range("D61").Value = 317.2
$ws.Range("E61").Value = 0.0473
$ws.Range("H61").Value = 11468.2385
$ws.Range("B62:C62").Value = "LTGILTCASE"
$ws.Range("D62").Value = 29.68
$ws.Range("E62").Value = 0
$ws.Range("H62").Value = 0
$ws.Range("B63:C63").Value = "REGAAL"
$ws.Range("D63").Value = 90.59999999999999
$ws.Range("E63").Value = 0.0663
$ws.Range("H63").Value = 930.0582000000001
$ws.Range("B64:C64").Value = "BLUESTONE"
$ws.Range("D64").Value = 571.75
$ws.Range("E64").Value = -0.522
$ws.Range("H64").Value = 8697.1384
$ws.Range("B65:C65").Value = "MOSILVER"
$ws.Range("D65").Value = 152.7
$ws.Range("E65").Value = -3.7504
$ws.Range("H65").Value = 0
$ws.Range("B66:C66").Value = "ALLTIME"
$ws.Range("D66").Value = 293.3
$ws.Range("E66").Value = -0.8284
$ws.Range("H66").Value = 1937.3925
$ws.Range("B67:C67").Value = "JSWCEMENT"
$ws.Range("D67").Value = 124.81
$ws.Range("E67").Value = -1.2423
$ws.Range("H67").Value = 17230.2061
$ws.Range("B68:C68").Value = "SBILIQETF"
$ws.Range("D68").Value = 1014.6
$ws.Range("E68").Value = 0.0138
$ws.Range("H68").Value = 0
$ws.Range("B69:C69").Value = "HILINFRA"
$ws.Range("D69").Value = 69.19
$ws.Range("E69").Value = 0.0145
$ws.Range("H69").Value = 0
$ws.Range("B70:C70").Value = "GROWWPOWER"
$ws.Range("D70").Value = 10.17
$ws.Range("E70").Value = 0.7929
$ws.Range("H70").Value = 0
$ws.Range("B71:C71").Value = "LOTUSDEV"
$ws.Range("D71").Value = 171.22
$ws.Range("E71").Value = -1.2287
$ws.Range("H71").Value = 8472.0255
$ws.Range("B72:C72").Value = "MBEL"
$ws.Range("D72").Value = 428.5
$ws.Range("E72").Value = -2.046
$ws.Range("H72").Value = 2499.9487
$ws.Range("B73:C73").Value = "LAXMIINDIA"
$ws.Range("D73").Value = 132.97
$ws.Range("E73").Value = -3.0265
$ws.Range("H73").Value = 716.6971
$ws.Range("B74:C74").Value = "CPPLUS"
$ws.Range("D74").Value = 1669.2
$ws.Range("E74").Value = 2.9798
$ws.Range("H74").Value = 19000.4952
$ws.Range("B75:C75").Value = "SHANTIGOLD"
$ws.Range("D75").Value = 226.87
$ws.Range("E75").Value = -0.3295
$ws.Range("H75").Value = 1641.0492
$ws.Range("B76:C76").Value = "MOGOLD"
$ws.Range("D76").Value = 121.35
$ws.Range("E76").Value = -2.491
$ws.Range("H76").Value = 0
